$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates: entrada/salida times and elapsed time changed
$ws.Range("C2").Value = "28/01/2025 09:29:40"
$ws.Range("D2").Value = "28/01/2025 09:36:20"
$ws.Range("G2").Value = "00:06:06"

# Row 3 updates: entrada/salida times, elapsed time, and total value changed
$ws.Range("C3").Value = "28/01/2025 09:35:45"
$ws.Range("D3").Value = "28/01/2025 09:36:49"
$ws.Range("G3").Value = "00:01:01"
$ws.Range("H3").Value = 4000

# Row 4 removed entirely (FEP0000003 / CCC01 entry deleted)
$ws.Rows.Item(4).Delete()
